$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19: 2012-09-10 (serial 41162), new string "Added LocalScan from Gpu gems 3"
$ws.Cells.Item(19, 1).Value = 41162
$ws.Cells.Item(19, 2).Value = "Added LocalScan from Gpu gems 3"

# Row 20: 2012-09-11 (serial 41163), new string "Further working on LocalScan"
$ws.Cells.Item(20, 1).Value = 41163
$ws.Cells.Item(20, 2).Value = "Further working on LocalScan"

# Row 21: 2012-09-12 (serial 41164), new string "LocalScan technically working now, allocation problems with sizes > 2MiB"
$ws.Cells.Item(21, 1).Value = 41164
$ws.Cells.Item(21, 2).Value = "LocalScan technically working now, allocation problems with sizes > 2MiB"

# Update the selected cell in the sheet view
$ws.Range("B24").Select()
